# Update crypto price (D) and 1h volume/change (E) columns with latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    # Force genuinely numeric-looking strings (e.g. "58.17") to be stored as
    # text -- matches the source data which keeps these as plain strings -- by
    # switching to a text number format just long enough to assign the value,
    # then reverting the style so no stray formatting is left behind.
    if ($Text -match '^[+-]?\d+(\.\d+)?$') {
        $Cell.NumberFormat = "@"
        $Cell.Value = $Text
        $Cell.Style = "Normal"
    } else {
        $Cell.Value = $Text
    }
}

Set-TextValue $ws.Range("D2") '70.320.77'
Set-TextValue $ws.Range("E2") '  +5.00%  '
Set-TextValue $ws.Range("D3") '3.611.26'
Set-TextValue $ws.Range("E3") '  +4.90%  '
Set-TextValue $ws.Range("E4") '  +0.04%  '
Set-TextValue $ws.Range("D5") '588.47'
Set-TextValue $ws.Range("E5") '  +3.16%  '
Set-TextValue $ws.Range("D6") '190.91'
Set-TextValue $ws.Range("E6") '  +3.22%  '
Set-TextValue $ws.Range("D7") '0.643'
Set-TextValue $ws.Range("E7") '  +1.53%  '
Set-TextValue $ws.Range("D8") '3.602.93'
Set-TextValue $ws.Range("E8") '  +4.79%  '
Set-TextValue $ws.Range("D9") '1.00'
Set-TextValue $ws.Range("E9") '  -0.02%  '
Set-TextValue $ws.Range("E10") '  +0.30%  '
Set-TextValue $ws.Range("D11") '0.661'
Set-TextValue $ws.Range("E11") '  +2.54%  '
Set-TextValue $ws.Range("D12") '58.17'
Set-TextValue $ws.Range("E12") '  +5.10%  '
Set-TextValue $ws.Range("E13") '  +3.54%  '
Set-TextValue $ws.Range("D14") '9.79'
Set-TextValue $ws.Range("E14") '  +4.19%  '
Set-TextValue $ws.Range("D15") '4.187.59'
Set-TextValue $ws.Range("E15") '  +4.99%  '
Set-TextValue $ws.Range("D16") '3.607.48'
Set-TextValue $ws.Range("E16") '  +4.96%  '
Set-TextValue $ws.Range("E17") '  +4.45%  '
Set-TextValue $ws.Range("D18") '70.198.37'
Set-TextValue $ws.Range("E18") '  +4.94%  '
Set-TextValue $ws.Range("E19") '  +3.89%  '
Set-TextValue $ws.Range("E20") '  +0.25%  '
Set-TextValue $ws.Range("E21") '  +3.98%  '
Set-TextValue $ws.Range("D22") '492.45'
Set-TextValue $ws.Range("E22") '  +2.88%  '
Set-TextValue $ws.Range("D23") '17.35'
Set-TextValue $ws.Range("E23") '  +13.88%  '
Set-TextValue $ws.Range("D24") '5.39'
Set-TextValue $ws.Range("E24") '  +8.23%  '
Set-TextValue $ws.Range("D25") '4.45'
Set-TextValue $ws.Range("E25") '  +6.33%  '
Set-TextValue $ws.Range("D26") '90.73'
Set-TextValue $ws.Range("E26") '  +1.17%  '
Set-TextValue $ws.Range("E27") '  +5.07%  '
Set-TextValue $ws.Range("D28") '11.12'
Set-TextValue $ws.Range("E28") '  +0.86%  '
Set-TextValue $ws.Range("D29") '9.47'
Set-TextValue $ws.Range("E29") '  +6.02%  '
Set-TextValue $ws.Range("D30") '32.42'
Set-TextValue $ws.Range("E30") '  +2.71%  '
Set-TextValue $ws.Range("D31") '7.54'
Set-TextValue $ws.Range("E31") '  +8.68%  '
Set-TextValue $ws.Range("D32") '627.13'
Set-TextValue $ws.Range("E32") '  +6.35%  '
Set-TextValue $ws.Range("D33") '12.25'
Set-TextValue $ws.Range("E33") '  +5.21%  '
Set-TextValue $ws.Range("D34") '0.117'
Set-TextValue $ws.Range("E34") '  +6.91%  '
Set-TextValue $ws.Range("D35") '65.38'
Set-TextValue $ws.Range("E35") '  +3.50%  '
Set-TextValue $ws.Range("E36") '  +4.15%  '
Set-TextValue $ws.Range("D37") '38.10'
Set-TextValue $ws.Range("E38") '  +3.65%  '
Set-TextValue $ws.Range("E39") '  +0.06%  '
Set-TextValue $ws.Range("E40") '  -1.20%  '
Set-TextValue $ws.Range("E41") '  -0.38%  '
Set-TextValue $ws.Range("D42") '3.300.58'
Set-TextValue $ws.Range("E42") '  +5.58%  '
Set-TextValue $ws.Range("E43") '  +5.77%  '
Set-TextValue $ws.Range("E44") '  +4.87%  '
Set-TextValue $ws.Range("E45") '  +1.52%  '
Set-TextValue $ws.Range("E46") '  +2.59%  '
Set-TextValue $ws.Range("E47") '  +1.81%  '
Set-TextValue $ws.Range("D48") '9.08'
Set-TextValue $ws.Range("E48") '  +4.37%  '
Set-TextValue $ws.Range("D49") '2.71'
Set-TextValue $ws.Range("E49") '  -3.37%  '
Set-TextValue $ws.Range("E50") '  +5.63%  '
Set-TextValue $ws.Range("D51") '0.999'
Set-TextValue $ws.Range("E51") '  -0.01%  '
